$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 192, pushing the
# existing rows 192-300 down to 193-301 (dimension grows from R300 to R301).
$ws.Rows.Item(192).Insert()

$ws.Range("A192").Value = 10
$ws.Range("B192").Value = "Vega Modelo de Temuco"
$ws.Range("C192").Value = "La Araucanía"
$ws.Range("D192").Value = 44518
$ws.Range("E192").Value = 9
$ws.Range("F192").Value = 100114014
$ws.Range("G192").Value = "Betarraga"
$ws.Range("H192").Value = "Sin especificar"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 115
$ws.Range("K192").Value = 9000
$ws.Range("L192").Value = 9000
$ws.Range("M192").Value = 9000
$ws.Range("N192").Value = "$/docena de paquetes"
$ws.Range("O192").Value = "Región del Maule"
$ws.Range("P192").Value = 750
$ws.Range("Q192").Value = 12
$ws.Range("R192").Value = "Hortaliza"
